# "se crea punto 42"
# Updates the SIPOT report row (record) to the next reporting quarter:
#  - new reporting period (start/end) dates
#  - new validation/update dates
#  - new explanatory note text about why there is no Consejo Consultivo
#  - matching cosmetic formatting (column width, row height, borders/alignment,
#    data validation range, active selection) that Excel applies when the
#    sheet is re-saved after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the explanatory note in K8 (shared string content) ---------
$ws.Range("K8").Value = "La Universidad Politécnica de Pachuca, no tiene Consejo Consultivo por que se encuentra en proceso de validación el nuevo Decreto de Creación."

# --- 2. Update the reporting period / validation dates (row 8) ------------
# B8/C8: period covered by this report (2023-04-01 .. 2023-06-30)
$ws.Range("B8").Value = 45017
$ws.Range("C8").Value = 45107
# I8/J8: validation date / update date (2023-07-10)
$ws.Range("I8").Value = 45117
$ws.Range("J8").Value = 45117

# --- 3. Widen column K so the longer note fits with fewer wrapped lines ---
$ws.Columns.Item(11).ColumnWidth = 67.16

# --- 4. Row 8 is now shorter since the text wraps into fewer lines --------
$ws.Rows.Item(8).RowHeight = 30

# --- 5. Refresh the cell formatting on row 8 -------------------------------
# D8:G8 lose their explicit left-alignment (back to default/general)
$ws.Range("D8:G8").Style = "Normal"
$ws.Range("D8:G8").Borders.LineStyle = 1

# K8 gets a full box border (instead of right/bottom only) and drops the
# explicit black Calibri font override, keeping justify + wrap text.
$ws.Range("K8").Style = "Normal"
$ws.Range("K8").Borders.LineStyle = 1
$ws.Range("K8").HorizontalAlignment = -4130
$ws.Range("K8").WrapText = $true

# --- 6. Data validation list now only spans the populated rows ------------
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D111").Validation.Add(3, 1, 1, "=Hidden_13")
$ws.Range("D8:D111").Validation.ShowInput = $false

# --- 7. Page setup / active selection (cosmetic, matches re-saved file) ---
$ws.PageSetup.Orientation = 1
$ws.Range("K14").Select()
